$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.01129674911499
$ws.Range("B1").Value = 2.128483295440674
$ws.Range("C1").Value = 5.807652950286865
$ws.Range("D1").Value = 1.028786659240723
$ws.Range("E1").Value = 1.100803017616272
